$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IESO project calibration: ON LDV ZEV market share_class_min drops slightly
# in the later years (2035-2050) from 100% to 99%.
$ws.Range("T3:W3").Value = 0.99
